# Added summary codes for trainings
#
# Rows 37-55 (HOR 3 "missing_trainings" entries, column E) were duplicating
# the Cluster value already present in column F. Remove the redundant
# "HOR 3" cell in column E for each of those rows, leaving the Cluster
# value (column F) and household_id (column G) untouched. Once the only
# references to the "HOR 3" shared string are gone, the workbook's shared
# string table drops that entry automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E37:E55").ClearContents()
